$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1212.3334
$ws.Range("I19").Value = 1425.2
$ws.Range("J19").Value = 1087.1177
$ws.Range("K19").Value = 1425.2
$ws.Range("L19").Value = 1087.1177
$ws.Range("M19").Value = -1250.2
$ws.Range("N19").Value = -1437.1177
$ws.Range("H33").Value = 35758084
$ws.Range("I33").Value = 52632652
$ws.Range("J33").Value = 133999.11
$ws.Range("K33").Value = 52632652
$ws.Range("L33").Value = 133999.11
$ws.Range("M33").Value = -52632423
$ws.Range("N33").Value = -134457.11
$ws.Range("H38").Value = 388.83334
$ws.Range("I38").Value = 99.89474
$ws.Range("J38").Value = 887.9091
$ws.Range("K38").Value = 299.68422
$ws.Range("L38").Value = 2663.7273
$ws.Range("M38").Value = 72.31578000000002
$ws.Range("N38").Value = -3407.7273
$ws.Range("H41").Value = 247.76471
$ws.Range("I41").Value = 103
$ws.Range("J41").Value = 454.57144
$ws.Range("K41").Value = 103
$ws.Range("L41").Value = 454.57144
$ws.Range("M41").Value = 337
$ws.Range("N41").Value = -1334.57144
$ws.Range("H53").Value = 127.84615
$ws.Range("I53").Value = 146.4
$ws.Range("J53").Value = 66
$ws.Range("K53").Value = 146.4
$ws.Range("L53").Value = 66
$ws.Range("M53").Value = 490.6
$ws.Range("N53").Value = -1340
$ws.Range("H98").Value = 57230.61
$ws.Range("I98").Value = 78350.08
$ws.Range("J98").Value = 2320
$ws.Range("K98").Value = 78350.08
$ws.Range("L98").Value = 2320
$ws.Range("M98").Value = -76852.08
$ws.Range("N98").Value = -5316
$ws.Range("H113").Value = 208417.72
$ws.Range("I113").Value = 319834.6
$ws.Range("J113").Value = 3176.1052
$ws.Range("K113").Value = 319834.6
$ws.Range("L113").Value = 3176.1052
$ws.Range("M113").Value = -316580.6
$ws.Range("N113").Value = -9684.1052
$ws.Range("H122").Value = 57230.61
$ws.Range("I122").Value = 78350.08
$ws.Range("J122").Value = 2320
$ws.Range("K122").Value = 235050.24
$ws.Range("L122").Value = 6960
$ws.Range("M122").Value = -232600.24
$ws.Range("N122").Value = -11860
$ws.Range("H132").Value = 1914.2778
$ws.Range("I132").Value = 872.2727
$ws.Range("K132").Value = 2616.8181
$ws.Range("M132").Value = -86.81809999999996
$ws.Range("H137").Value = 1133.3438
$ws.Range("I137").Value = 1103.9131
$ws.Range("J137").Value = 1208.5555
$ws.Range("K137").Value = 3311.7393
$ws.Range("L137").Value = 3625.6665
$ws.Range("M137").Value = -761.7393000000002
$ws.Range("N137").Value = -8725.666499999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9981.48
$ws.Range("I32").Value = 4002.2632
$ws.Range("J32").Value = 28915.666
$ws.Range("K32").Value = 4002.2632
$ws.Range("L32").Value = 28915.666
$ws.Range("M32").Value = -3715.2632
$ws.Range("N32").Value = -29489.666
$ws.Range("H61").Value = 3166.6667
$ws.Range("I61").Value = 3250
$ws.Range("J61").Value = 3125
$ws.Range("K61").Value = 3250
$ws.Range("L61").Value = 3125
$ws.Range("M61").Value = -3038
$ws.Range("N61").Value = -3549
$ws.Range("H136").Value = 3166.6667
$ws.Range("I136").Value = 3250
$ws.Range("J136").Value = 3125
$ws.Range("K136").Value = 9750
$ws.Range("L136").Value = 9375
$ws.Range("M136").Value = -7200
$ws.Range("N136").Value = -14475

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 855.44446
$ws.Range("I94").Value = 712.375
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 712.375
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -261.375
$ws.Range("N94").Value = -2902

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 16584.166
$ws.Range("I2").Value = 8666.666999999999
$ws.Range("J2").Value = 24501.666
$ws.Range("K2").Value = 8666.666999999999
$ws.Range("L2").Value = 24501.666
$ws.Range("M2").Value = -8553.666999999999
$ws.Range("N2").Value = -24727.666
$ws.Range("H3").Value = 500000
$ws.Range("I3").Value = 500000
$ws.Range("K3").Value = 500000
$ws.Range("M3").Value = -499887
$ws.Range("H4").Value = 5225500
$ws.Range("J4").Value = 5225500
$ws.Range("L4").Value = 5225500
$ws.Range("N4").Value = -5225724
$ws.Range("H16").Value = 6422.5
$ws.Range("I16").Value = 4716
$ws.Range("J16").Value = 9266.666999999999
$ws.Range("K16").Value = 4716
$ws.Range("L16").Value = 9266.666999999999
$ws.Range("M16").Value = -4429
$ws.Range("N16").Value = -9840.666999999999
$ws.Range("H31").Value = 2817.45
$ws.Range("I31").Value = 2460.875
$ws.Range("K31").Value = 2460.875
$ws.Range("M31").Value = -2165.875
$ws.Range("H34").Value = 2817.45
$ws.Range("I34").Value = 2460.875
$ws.Range("K34").Value = 2460.875
$ws.Range("M34").Value = -2258.875
$ws.Range("H75").Value = 48260
$ws.Range("J75").Value = 48260
$ws.Range("L75").Value = 48260
$ws.Range("N75").Value = -50256
$ws.Range("H78").Value = 48260
$ws.Range("J78").Value = 48260
$ws.Range("L78").Value = 144780
$ws.Range("N78").Value = -154764
$ws.Range("H113").Value = 6422.5
$ws.Range("I113").Value = 4716
$ws.Range("J113").Value = 9266.666999999999
$ws.Range("K113").Value = 4716
$ws.Range("L113").Value = 9266.666999999999
$ws.Range("M113").Value = -2546
$ws.Range("N113").Value = -13606.667
$ws.Range("H132").Value = 2788.6
$ws.Range("I132").Value = 977.8
$ws.Range("J132").Value = 4599.4
$ws.Range("K132").Value = 2933.4
$ws.Range("L132").Value = 13798.2
$ws.Range("M132").Value = -403.3999999999996
$ws.Range("N132").Value = -18858.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42736.293
$ws.Range("I4").Value = 222.33333
$ws.Range("K4").Value = 666.99999
$ws.Range("M4").Value = -554.99999
$ws.Range("H5").Value = 594.75
$ws.Range("I5").Value = 499.5
$ws.Range("K5").Value = 1498.5
$ws.Range("M5").Value = -1386.5
$ws.Range("H6").Value = 136.61111
$ws.Range("I6").Value = 78.6875
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 236.0625
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -123.0625
$ws.Range("N6").Value = -2026
$ws.Range("H113").Value = 4762321.5
$ws.Range("I113").Value = 356.25
$ws.Range("J113").Value = 10204567
$ws.Range("K113").Value = 1068.75
$ws.Range("L113").Value = 30613701
$ws.Range("M113").Value = 1101.25
$ws.Range("N113").Value = -30618041
$ws.Range("H122").Value = 2947853
$ws.Range("I122").Value = 777.6667
$ws.Range("J122").Value = 6263313
$ws.Range("K122").Value = 6999.0003
$ws.Range("L122").Value = 56369817
$ws.Range("M122").Value = -4549.0003
$ws.Range("N122").Value = -56374717
$ws.Range("H135").Value = 594.75
$ws.Range("I135").Value = 499.5
$ws.Range("K135").Value = 4495.5
$ws.Range("M135").Value = -1960.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.75
$ws.Range("I2").Value = 28.666666
$ws.Range("J2").Value = 125
$ws.Range("K2").Value = 28.666666
$ws.Range("L2").Value = 125
$ws.Range("M2").Value = 84.33333400000001
$ws.Range("N2").Value = -351
$ws.Range("H3").Value = 3203
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 70004
$ws.Range("J4").Value = 70004
$ws.Range("L4").Value = 70004
$ws.Range("N4").Value = -70228
$ws.Range("H80").Value = 2928.6
$ws.Range("I80").Value = 2760
$ws.Range("J80").Value = 3097.2
$ws.Range("K80").Value = 2760
$ws.Range("L80").Value = 3097.2
$ws.Range("M80").Value = -1762
$ws.Range("N80").Value = -5093.2
$ws.Range("H83").Value = 2928.6
$ws.Range("I83").Value = 2760
$ws.Range("J83").Value = 3097.2
$ws.Range("K83").Value = 13800
$ws.Range("L83").Value = 15486
$ws.Range("M83").Value = -8808
$ws.Range("N83").Value = -25470

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1125889
$ws.Range("J2").Value = 16625
$ws.Range("L2").Value = 16625
$ws.Range("N2").Value = -16849
$ws.Range("H132").Value = 2863.6667
$ws.Range("I132").Value = 2224.889
$ws.Range("J132").Value = 4013.4666
$ws.Range("K132").Value = 6674.667
$ws.Range("L132").Value = 12040.3998
$ws.Range("M132").Value = -4144.667
$ws.Range("N132").Value = -17100.3998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1623.44
$ws.Range("I132").Value = 910.7368
$ws.Range("J132").Value = 2060.258
$ws.Range("K132").Value = 2732.2104
$ws.Range("L132").Value = 6180.773999999999
$ws.Range("M132").Value = -202.2103999999999
$ws.Range("N132").Value = -11240.774
$ws.Range("H136").Value = 1325.4348
$ws.Range("I136").Value = 969.41174
$ws.Range("J136").Value = 2334.1667
$ws.Range("K136").Value = 2908.23522
$ws.Range("L136").Value = 7002.500100000001
$ws.Range("M136").Value = -358.23522
$ws.Range("N136").Value = -12102.5001
